$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table 1")

# Swap the header labels in B1 and C1 (V+ and V- order changed)
$ws.Range("B1").Value = "V-"
$ws.Range("C1").Value = "V+"

# Update the active selection to C1 (as reflected in the saved file)
$ws.Range("C1").Select()
